# -----------------------------------------------------------------------
# Applies the citation clean-up edit described by the commit:
#   - "Fuks  and" keeps its text but is split into two runs
#   - curly quotes removed around the paper title, and the stray period
#     right before the closing quote is dropped
#   - the "_GoBack" bookmark is relocated to sit right after the (now
#     unquoted) paper title instead of at the end of the Results section
#   - curly quotes removed around the method names 'rpart' and 'rf'
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

$lq = [char]0x201C   # U+201C LEFT DOUBLE QUOTATION MARK  "
$rq = [char]0x201D   # U+201D RIGHT DOUBLE QUOTATION MARK "
$oq = [char]0x2018   # U+2018 LEFT SINGLE QUOTATION MARK  '
$cq = [char]0x2019   # U+2019 RIGHT SINGLE QUOTATION MARK '

# Insert+remove a same-named bookmark at a position to force the engine
# to split the run that currently spans that position into two runs,
# without leaving any trace behind.
function Split-At($pos, $name) {
    $tmp = $d.Range($pos, $pos)
    $d.Bookmarks.Add($name, $tmp)
    $d.Bookmarks($name).Delete()
}

# ------------------------------------------------------------------
# 1) Split the "  and" run into "  " + "and" (no visible text change).
# ------------------------------------------------------------------
$rSplit = $d.Content
$foundSplit = $rSplit.Find.Execute("  and", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundSplit) {
    Split-At ($rSplit.Start + 2) "splitAnd"
}

# ------------------------------------------------------------------
# 2) Remove the curly quotes around the citation title, and drop the
#    trailing period that sat just inside the closing quote:
#       ... presented in "Qualitative ... Exercises."; Proceedings ...
#    becomes
#       ... presented in Qualitative ... Exercises; Proceedings ...
# ------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("presented in " + $lq, $true, $false, $false, $false, $false, $true, 1, $false, "presented in ", 2)

$r2 = $d.Content
$found2 = $r2.Find.Execute("Exercises." + $rq + "; Proceedings", $true, $false, $false, $false, $false, $true, 1, $false, "Exercises; Proceedings", 2)

# The edit above merges the trailing URL run back into the sentence run;
# split it back out so the hyperlink text keeps its own run, as before.
$r2b = $d.Content
$found2b = $r2b.Find.Execute("downloaded from ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2b) {
    Split-At $r2b.End "splitUrl"
}

# ------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark so it now sits right after
#    "Weight Lifting Exercises;" instead of at the end of the Random
#    Forest paragraph (Bookmarks.Add silently replaces the bookmark
#    that previously held this name).
# ------------------------------------------------------------------
$r3 = $d.Content
$found3 = $r3.Find.Execute("Weight Lifting Exercises;", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $bkRange = $d.Range($r3.End, $r3.End)
    $d.Bookmarks.Add("_GoBack", $bkRange)
}

# ------------------------------------------------------------------
# 4) Remove the curly single quotes around 'rpart'.
# ------------------------------------------------------------------
$r4 = $d.Content
$found4 = $r4.Find.Execute("with method " + $oq, $true, $false, $false, $false, $false, $true, 1, $false, "with method ", 2)

$r5 = $d.Content
$found5 = $r5.Find.Execute($cq + ", due to its ease", $true, $false, $false, $false, $false, $true, 1, $false, ", due to its ease", 2)

# ------------------------------------------------------------------
# 5) Remove the curly single quotes around 'rf'.
# ------------------------------------------------------------------
$r6 = $d.Content
$found6 = $r6.Find.Execute("using method " + $oq, $true, $false, $false, $false, $false, $true, 1, $false, "using method ", 2)

$r7 = $d.Content
$found7 = $r7.Find.Execute($cq + ". The Random", $true, $false, $false, $false, $false, $true, 1, $false, ". The Random", 2)

# The edit above merges the rest of that paragraph into a single run;
# split it back apart at the original run boundaries.
if ($found7) {
    $scoped = $d.Range($r7.End, $d.Content.End)
    $find8 = $scoped.Duplicate
    $f8 = $find8.Find.Execute("orest Model performed very ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($f8) {
        Split-At $find8.Start "s1"
        Split-At $find8.End "s2"

        $find9 = $d.Range($find8.End, $d.Content.End).Duplicate
        $f9 = $find9.Find.Execute("well ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($f9) { Split-At $find9.End "s3" }

        $find10 = $d.Range($find9.End, $d.Content.End).Duplicate
        $f10 = $find10.Find.Execute("with ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($f10) { Split-At $find10.End "s4" }

        $find11 = $d.Range($find10.End, $d.Content.End).Duplicate
        $f11 = $find11.Find.Execute("reported ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($f11) { Split-At $find11.End "s5" }

        $find12 = $d.Range($find11.End, $d.Content.End).Duplicate
        $f12 = $find12.Find.Execute("accuracy of ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($f12) { Split-At $find12.End "s6" }

        $find13 = $d.Range($find12.End, $d.Content.End).Duplicate
        $f13 = $find13.Find.Execute("98.2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($f13) { Split-At $find13.End "s7" }
    }
}

Write-Output "split=$foundSplit quote1=$found1 quote2=$found2 splitUrl=$found2b bookmark=$found3 rpart1=$found4 rpart2=$found5 rf1=$found6 rf2=$found7"
